# Auto-context: $ppt.ActivePresentation is already open as $p
$p = $ppt.ActivePresentation

# ---- Slide 4: merge the two trailing runs in the last bullet paragraph ----
$s4 = $p.Slides.Item(4)
$body4 = $s4.Shapes.Item(2).TextFrame.TextRange
$lastPara = $body4.Paragraphs($body4.Paragraphs().Count, 1)
$lastPara.Text = "ZZZ_TEMP_PLACEHOLDER_ZZZ"
$lastPara = $body4.Paragraphs($body4.Paragraphs().Count, 1)
$lastPara.Runs(1).Text = "State: closed-source commercial library, then declassified (circa 2016) and placed to the open-source. "

# ---- Slide 5 ----
$s = $p.Slides.Add(5, 2)
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Test code"
$bodyRange = $s.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "Quite a part of pbl_met is “test code”, aimed at harnessing procedure functionalities, and “proving” they are correct.`rThis is unlike the legacy PBL_MET, for which test code was not released.`rTest code, collected under directory “/test”, also provides examples on using individual routines."
$bodyRange.Characters(17, 7).Font.Italic = $true

# ---- Slide 6 ----
$s = $p.Slides.Add(6, 2)
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Why fortran?"
$bodyRange = $s.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "To date, much scientific code is written using Python and R, so why indulging on Fortran?`rOur main reasons are, small footprint and run-time efficiency.`rLikely use cases:`rReal-time met processors on small scale embedded systems.`rInclusion as component in meteorological and atmospheric pollutant dispersion models.`rBesides, Fortran is still used and well known among the geophysicists community."
$bodyRange.Paragraphs(4, 1).IndentLevel = 2
$bodyRange.Paragraphs(5, 1).IndentLevel = 2

# ---- Slide 7 ----
$s = $p.Slides.Add(7, 2)
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Fortran specific advantages"
$bodyRange = $s.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "In modern Fortran (i.e. Fortran 2003, 2008, 2015) extensive support is provided for current programming paradigms, e.g. object orientation and large-scale programming.`rModern Fortran natively supports parallelism, both fine-grained (e.g. SIMD) and coarse-grained (e.g. multi-core).`rNot directly used in pbl_met, but allowing easy inclusion in parallel code.`rVery important (most, possibly): modern Fortran allows writing readable and understandable code."
$bodyRange.Paragraphs(3, 1).IndentLevel = 2
$bodyRange.Characters(4, 6).Font.Italic = $true
$bodyRange.Characters(304, 7).Font.Italic = $true
$bodyRange.Characters(422, 8).Font.Italic = $true
$bodyRange.Characters(435, 14).Font.Italic = $true

# ---- Slide 8 ----
$s = $p.Slides.Add(8, 2)
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Code readability"
$bodyRange = $s.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "The pbl_met places a very high value in readable and understandable code.`rOur phylosophy departs somewhat from traditional open-source: we firmly believe there is no logical reason people, in front of a “public” project, should “sink or swim”.`rOur purpose is not selecting tough contributors, but rather making life the least miserable possible to prospective users, who are not necessarily dedicated amateurs or professional programmers.`rBecause of this, pbl_met is “written for people” instead of “for the machine”."
$bodyRange.Characters(5, 7).Font.Italic = $true
$bodyRange.Characters(41, 8).Font.Italic = $true
$bodyRange.Characters(54, 14).Font.Italic = $true
$bodyRange.Characters(457, 7).Font.Italic = $true

